$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.515.01"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.627.06"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "212.91"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "18.82"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "1.853.49"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.623.84"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "65.05"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "26.528.36"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "214.51"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").Value = "9.29"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +9.18%  "
$ws.Range("D25").Value = "148.35"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "1.240.52"
$ws.Range("E34").Value = "  +6.16%  "
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("E37").Value = "  +3.96%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "1.764.83"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "92.96"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "1.58"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "54.86"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.406"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  -0.84%  "
